$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: add the NeoPixel / "LED strip" line item quantity, unit cost and link ---
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 23.2
$ws.Range("E9").Value = "https://onecall.farnell.com/mcm/25-7520/30-rgb-led-addressable-strip-1/dp/2801872?st=individually%20addressable%20rgb%20led%20strip"

# --- Attach a hyperlink to the existing MCP3301 link text in E3 (it was plain text before) ---
$e3Text = $ws.Range("E3").Text
$ws.Hyperlinks.Add($ws.Range("E3"), "https://onecall.farnell.com/microchip/mcp3301-ci-p/ic-13bit-adc-1ch-dip8-3301/dp/1332099", "", "", $e3Text)
# Hyperlinks.Add re-styles the cell (bold/underline/color) as a side-effect - restore the original look.
$ws.Range("E3").Font.Bold = $false
$ws.Range("E3").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleNone
$ws.Range("E3").Font.Color = 0

# --- Row height tweaks ---
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(9).RowHeight = 13.8

# --- Move the active selection ---
$ws.Range("A11").Select()
